$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(4, 6).Value = 1778
$ws.Cells.Item(4, 7).Value = 70
$ws.Cells.Item(6, 6).Value = 627
$ws.Cells.Item(7, 6).Value = 1153
$ws.Cells.Item(8, 6).Value = 1585
$ws.Cells.Item(9, 6).Value = 175
$ws.Cells.Item(10, 6).Value = 175
$ws.Cells.Item(12, 6).Value = 1514
$ws.Cells.Item(13, 6).Value = 3135
$ws.Cells.Item(14, 6).Value = 676
$ws.Cells.Item(15, 6).Value = 1830
$ws.Cells.Item(16, 6).Value = 1826
$ws.Cells.Item(17, 6).Value = 888
$ws.Cells.Item(18, 6).Value = 301
$ws.Cells.Item(20, 6).Value = 1510
$ws.Cells.Item(21, 6).Value = 307
$ws.Cells.Item(23, 6).Value = 26
$ws.Cells.Item(24, 6).Value = 1295
$ws.Cells.Item(25, 6).Value = 421
$ws.Cells.Item(26, 6).Value = 506
$ws.Cells.Item(27, 6).Value = 187
$ws.Cells.Item(28, 6).Value = 6792
$ws.Cells.Item(29, 6).Value = 5401
$ws.Cells.Item(30, 6).Value = 772
$ws.Cells.Item(32, 6).Value = 1710
$ws.Cells.Item(33, 6).Value = 92
$ws.Cells.Item(34, 6).Value = 235

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 6).Value = 1

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 50

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(3, 6).Value = 50
$ws.Cells.Item(7, 6).Value = 1778
$ws.Cells.Item(7, 7).Value = 70
$ws.Cells.Item(9, 6).Value = 627
$ws.Cells.Item(10, 6).Value = 1153
$ws.Cells.Item(11, 6).Value = 1585
$ws.Cells.Item(12, 6).Value = 175
$ws.Cells.Item(13, 6).Value = 175
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(16, 6).Value = 1514
$ws.Cells.Item(17, 6).Value = 3135
$ws.Cells.Item(18, 6).Value = 676
$ws.Cells.Item(19, 6).Value = 1830
$ws.Cells.Item(20, 6).Value = 1826
$ws.Cells.Item(21, 6).Value = 888
$ws.Cells.Item(22, 6).Value = 301
$ws.Cells.Item(24, 6).Value = 1510
$ws.Cells.Item(25, 6).Value = 307
$ws.Cells.Item(28, 6).Value = 26
$ws.Cells.Item(30, 6).Value = 1295
$ws.Cells.Item(31, 6).Value = 421
$ws.Cells.Item(32, 6).Value = 506
$ws.Cells.Item(33, 6).Value = 187
$ws.Cells.Item(34, 6).Value = 6792
$ws.Cells.Item(35, 6).Value = 5401
$ws.Cells.Item(36, 6).Value = 772
$ws.Cells.Item(38, 6).Value = 1710
$ws.Cells.Item(41, 6).Value = 92
$ws.Cells.Item(42, 6).Value = 235
